# Update "SAM 2014 beta feedback" tracker: append a new feedback row (row 55)
# describing Jeff Cook's email about wind-farm turbine layout, plus the
# follow-up note, and move the active selection to the next empty row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- New row 55 --------------------------------------------------------
# A: Date received
$ws.Range("A55").Value = 41929
$ws.Range("A55").NumberFormat = "m/d/yyyy"

# B: Source
$ws.Range("B55").Value = "Email from SAM"

# C: Contact
$ws.Range("C55").Value = "Cook, Jeff <Jeff.Cook@nrel.gov>"

# D: Description (wrapped)
$ws.Range("D55").Value = "I am a little curious as to why you can only adjust the rows and number of turbines per row, to get the total number of turbines to change. To me it would be more intuitive if you could change the number of turbines in the first cell, and have that adjust your rows and number of turbines per row by default. "
$ws.Range("D55").WrapText = $true

# E: Status / action taken (wrapped)
$ws.Range("E55").Value = "Replied with cc to Janine. Should be addressed by new Wind Farm layout option"
$ws.Range("E55").WrapText = $true

# F: Last reply date
$ws.Range("F55").Value = 41929
$ws.Range("F55").NumberFormat = "m/d/yyyy"

# Row height to match the wrapped, multi-line content of the other rows
$ws.Range("A55:F55").RowHeight = 75

# --- View state ----------------------------------------------------------
# Move the active cell/selection down to the next blank row, as happens
# after finishing data entry on row 55.
$ws.Range("E56").Select()
